$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.134879666666667
$ws.Range("H2").Value = 21.404639
$ws.Range("I2").Value = 0.07716103050836744
$ws.Range("J2").Value = 0.07716103050836744
$ws.Range("M2").Value = 15.959554
$ws.Range("N2").Value = 31.919108
$ws.Range("O2").Value = 0.1230717359814782
$ws.Range("P2").Value = 0.09235606860026654
$ws.Range("Q2").Value = 113.8694973236687
$ws.Range("R2").Value = 683.2169839420119
$ws.Range("S2").Value = 0.009496341974784585
$ws.Range("T2").Value = 0.007126289426898043
$ws.Range("G3").Value = 7.134879666666667
$ws.Range("H3").Value = 21.404639
$ws.Range("I3").Value = 0.07716103050836744
$ws.Range("J3").Value = 0.07716103050836744
$ws.Range("O3").Value = 0.6576648948789622
$ws.Range("P3").Value = 0.7402919565126984
$ws.Range("Q3").Value = 608.4904091916045
$ws.Range("R3").Value = 5476.413682724441
$ws.Range("S3").Value = 0.05074610101803786
$ws.Range("T3").Value = 0.05712169024157534
$ws.Range("G4").Value = 7.134879666666667
$ws.Range("H4").Value = 21.404639
$ws.Range("I4").Value = 0.07716103050836744
$ws.Range("J4").Value = 0.07716103050836744
$ws.Range("M4").Value = 0.4926256666666666
$ws.Range("N4").Value = 1.477877
$ws.Range("O4").Value = 0.00379887157158024
$ws.Range("P4").Value = 0.004276150498778228
$ws.Range("Q4").Value = 3.514824852378111
$ws.Range("R4").Value = 31.633423671403
$ws.Range("S4").Value = 0.0002931248452320726
$ws.Range("T4").Value = 0.0003299521790945975
$ws.Range("G5").Value = 7.134879666666667
$ws.Range("H5").Value = 21.404639
$ws.Range("I5").Value = 0.07716103050836744
$ws.Range("J5").Value = 0.07716103050836744
$ws.Range("M5").Value = 27.4617595
$ws.Range("N5").Value = 54.923519
$ws.Range("O5").Value = 0.2117707308594496
$ws.Range("P5").Value = 0.1589179838149626
$ws.Range("Q5").Value = 195.9363494674402
$ws.Range("R5").Value = 1175.618096804641
$ws.Range("S5").Value = 0.01634044782462526
$ws.Range("T5").Value = 0.01226227539747457
$ws.Range("G6").Value = 7.134879666666667
$ws.Range("H6").Value = 21.404639
$ws.Range("I6").Value = 0.07716103050836744
$ws.Range("J6").Value = 0.07716103050836744
$ws.Range("M6").Value = 0.363283
$ws.Range("N6").Value = 1.089849
$ws.Range("O6").Value = 0.002801448553171308
$ws.Range("P6").Value = 0.003153414218465375
$ws.Range("Q6").Value = 2.591980489945666
$ws.Range("R6").Value = 23.327824409511
$ws.Range("S6").Value = 0.0002161626572788731
$ws.Range("T6").Value = 0.0002433206907165264
$ws.Range("G7").Value = 7.134879666666667
$ws.Range("H7").Value = 21.404639
$ws.Range("I7").Value = 0.07716103050836744
$ws.Range("J7").Value = 0.07716103050836744
$ws.Range("M7").Value = 0.115713
$ws.Range("N7").Value = 0.347139
$ws.Range("O7").Value = 0.0008923181553585264
$ws.Range("P7").Value = 0.001004426354828836
$ws.Range("Q7").Value = 0.8255983308689999
$ws.Range("R7").Value = 7.430384977820999
$ws.Range("S7").Value = 0.00006885218840878941
$ws.Range("T7").Value = 0.00007750257260835608
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 81.06813666666666
$ws.Range("H8").Value = 243.20441
$ws.Range("I8").Value = 0.8767212985829614
$ws.Range("J8").Value = 0.8767212985829616
$ws.Range("M8").Value = 15.959554
$ws.Range("N8").Value = 31.919108
$ws.Range("O8").Value = 0.1230717359814782
$ws.Range("P8").Value = 0.09235606860026654
$ws.Range("Q8").Value = 1293.811304811047
$ws.Range("R8").Value = 7762.86782886628
$ws.Range("S8").Value = 0.107899612188541
$ws.Range("T8").Value = 0.08097053239524277
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 81.06813666666666
$ws.Range("H9").Value = 243.20441
$ws.Range("I9").Value = 0.8767212985829614
$ws.Range("J9").Value = 0.8767212985829616
$ws.Range("O9").Value = 0.6576648948789622
$ws.Range("P9").Value = 0.7402919565126984
$ws.Range("Q9").Value = 6913.807374097865
$ws.Range("R9").Value = 62224.26636688079
$ws.Range("S9").Value = 0.5765888206707106
$ws.Range("T9").Value = 0.6490297254443342
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 81.06813666666666
$ws.Range("H10").Value = 243.20441
$ws.Range("I10").Value = 0.8767212985829614
$ws.Range("J10").Value = 0.8767212985829616
$ws.Range("M10").Value = 0.4926256666666666
$ws.Range("N10").Value = 1.477877
$ws.Range("O10").Value = 0.00379887157158024
$ws.Range("P10").Value = 0.004276150498778228
$ws.Range("Q10").Value = 39.93624487084111
$ws.Range("R10").Value = 359.42620383757
$ws.Range("S10").Value = 0.003330551617385723
$ws.Range("T10").Value = 0.003748992218225027
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 81.06813666666666
$ws.Range("H11").Value = 243.20441
$ws.Range("I11").Value = 0.8767212985829614
$ws.Range("J11").Value = 0.8767212985829616
$ws.Range("M11").Value = 27.4617595
$ws.Range("N11").Value = 54.923519
$ws.Range("O11").Value = 0.2117707308594496
$ws.Range("P11").Value = 0.1589179838149626
$ws.Range("Q11").Value = 2226.273672253131
$ws.Range("R11").Value = 13357.64203351879
$ws.Range("S11").Value = 0.1856639101609594
$ws.Range("T11").Value = 0.1393267811384401
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 81.06813666666666
$ws.Range("H12").Value = 243.20441
$ws.Range("I12").Value = 0.8767212985829614
$ws.Range("J12").Value = 0.8767212985829616
$ws.Range("M12").Value = 0.363283
$ws.Range("N12").Value = 1.089849
$ws.Range("O12").Value = 0.002801448553171308
$ws.Range("P12").Value = 0.003153414218465375
$ws.Range("Q12").Value = 29.45067589267666
$ws.Range("R12").Value = 265.05608303409
$ws.Range("S12").Value = 0.002456089613449707
$ws.Range("T12").Value = 0.002764665408582938
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 81.06813666666666
$ws.Range("H13").Value = 243.20441
$ws.Range("I13").Value = 0.8767212985829614
$ws.Range("J13").Value = 0.8767212985829616
$ws.Range("M13").Value = 0.115713
$ws.Range("N13").Value = 0.347139
$ws.Range("O13").Value = 0.0008923181553585264
$ws.Range("P13").Value = 0.001004426354828836
$ws.Range("Q13").Value = 9.380637298109999
$ws.Range("R13").Value = 84.42573568298999
$ws.Range("S13").Value = 0.00078231433191508
$ws.Range("T13").Value = 0.0008806019781364873
$ws.Range("G14").Value = 4.264381
$ws.Range("H14").Value = 12.793143
$ws.Range("I14").Value = 0.04611767090867112
$ws.Range("J14").Value = 0.04611767090867113
$ws.Range("M14").Value = 15.959554
$ws.Range("N14").Value = 31.919108
$ws.Range("O14").Value = 0.1230717359814782
$ws.Range("P14").Value = 0.09235606860026654
$ws.Range("Q14").Value = 68.05761884607399
$ws.Range("R14").Value = 408.345713076444
$ws.Range("S14").Value = 0.005675781818152672
$ws.Range("T14").Value = 0.004259246778125747
$ws.Range("G15").Value = 4.264381
$ws.Range("H15").Value = 12.793143
$ws.Range("I15").Value = 0.04611767090867112
$ws.Range("J15").Value = 0.04611767090867113
$ws.Range("O15").Value = 0.6576648948789622
$ws.Range("P15").Value = 0.7402919565126984
$ws.Range("Q15").Value = 363.6830697736464
$ws.Range("R15").Value = 3273.147627962817
$ws.Range("S15").Value = 0.03032997319021377
$ws.Range("T15").Value = 0.03414054082678891
$ws.Range("G16").Value = 4.264381
$ws.Range("H16").Value = 12.793143
$ws.Range("I16").Value = 0.04611767090867112
$ws.Range("J16").Value = 0.04611767090867113
$ws.Range("M16").Value = 0.4926256666666666
$ws.Range("N16").Value = 1.477877
$ws.Range("O16").Value = 0.00379887157158024
$ws.Range("P16").Value = 0.004276150498778228
$ws.Range("Q16").Value = 2.100743533045667
$ws.Range("R16").Value = 18.906691797411
$ws.Range("S16").Value = 0.0001751951089624438
$ws.Range("T16").Value = 0.0001972061014586042
$ws.Range("G17").Value = 4.264381
$ws.Range("H17").Value = 12.793143
$ws.Range("I17").Value = 0.04611767090867112
$ws.Range("J17").Value = 0.04611767090867113
$ws.Range("M17").Value = 27.4617595
$ws.Range("N17").Value = 54.923519
$ws.Range("O17").Value = 0.2117707308594496
$ws.Range("P17").Value = 0.1589179838149626
$ws.Range("Q17").Value = 117.1074054383695
$ws.Range("R17").Value = 702.644432630217
$ws.Range("S17").Value = 0.00976637287386486
$ws.Range("T17").Value = 0.007328927279047971
$ws.Range("G18").Value = 4.264381
$ws.Range("H18").Value = 12.793143
$ws.Range("I18").Value = 0.04611767090867112
$ws.Range("J18").Value = 0.04611767090867113
$ws.Range("M18").Value = 0.363283
$ws.Range("N18").Value = 1.089849
$ws.Range("O18").Value = 0.002801448553171308
$ws.Range("P18").Value = 0.003153414218465375
$ws.Range("Q18").Value = 1.549177122823
$ws.Range("R18").Value = 13.942594105407
$ws.Range("S18").Value = 0.0001291962824427272
$ws.Range("T18").Value = 0.0001454281191659105
$ws.Range("G19").Value = 4.264381
$ws.Range("H19").Value = 12.793143
$ws.Range("I19").Value = 0.04611767090867112
$ws.Range("J19").Value = 0.04611767090867113
$ws.Range("M19").Value = 0.115713
$ws.Range("N19").Value = 0.347139
$ws.Range("O19").Value = 0.0008923181553585264
$ws.Range("P19").Value = 0.001004426354828836
$ws.Range("Q19").Value = 0.493444318653
$ws.Range("R19").Value = 4.440998867877
$ws.Range("S19").Value = 0.00004115163503465699
$ws.Range("T19").Value = 0.00004632180408399238
